$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column before column N (shifts N:P -> O:Q, R -> S, etc.)
$ws.Columns("N").Insert()

# The newly inserted column keeps a distinct (non bestFit) width, matching
# column M's width of 10.7109375 characters.
$ws.Columns("N").ColumnWidth = 9.83

# Restore the active selection shown in the saved workbook.
[void]$ws.Range("K18").Select()
